# fix caching in vignette
# - Update/append the "Hospitilization Data" sheet with refreshed data
# - Switch the active/selected tab from "Parameters with Distributions"
#   to "Hospitilization Data"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hospitilization Data")

# Make this the active sheet (moves tabSelected / sets workbook activeTab)
[void]$ws.Activate()

# --- Precision fix for two existing values ---
$ws.Range("C2").Value2 = 62.6
$ws.Range("C4").Value2 = 62.2

# --- Refresh the daily data table (rows 9-23 updated, rows 24-38 appended) ---
$data = @(
    @(9,  43922, 79, 103.9),
    @(10, 43923, 76, 92.2),
    @(11, 43926, 82, 118.6),
    @(12, 43927, 93, 128.69999999999999),
    @(13, 43928, 91, 107.8),
    @(14, 43929, 89, 105.2),
    @(15, 43930, 86, 98),
    @(16, 43931, 83, 94.7),
    @(17, 43932, 94, 104.5),
    @(18, 43933, 90, 102.6),
    @(19, 43934, 88, 100.9),
    @(20, 43935, 88, 113.8),
    @(21, 43936, 89, 100.4),
    @(22, 43937, 84, 94.5),
    @(23, 43938, 77, 85.1),
    @(24, 43939, 78, 88.8),
    @(25, 43940, 84, 93.6),
    @(26, 43941, 81, 94.8),
    @(27, 43942, 80, 88.4),
    @(28, 43943, 82, 90.4),
    @(29, 43944, 77, 84.5),
    @(30, 43945, 78, 84.6),
    @(31, 43946, 85, 97.6),
    @(32, 43947, 87, 111.9),
    @(33, 43948, 89, 100.7),
    @(34, 43949, 84, 92.4),
    @(35, 43950, 90, 96.6),
    @(36, 43951, 90, 96),
    @(37, 43952, 87, 95.7),
    @(38, 43953, 90, 96)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value2 = $row[1]
    $ws.Cells.Item($r, 2).Value2 = $row[2]
    $ws.Cells.Item($r, 3).Value2 = $row[3]
}

# --- Selection on the now-active sheet ---
[void]$ws.Range("G31").Select()

Write-Output "done"
